{"js": "// Update the p-values in the Mantel correlogram table (x-axis label /\n// distance-class rows). Four cells in the \"p\" column change:\n//   0.257 -> 0.26\n//   0.513 -> 0.519\n//   0.77  -> 0.779\n//   0.279 -> 0.261\n\nasync function replaceCellValue(oldText, newText) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  // Replace only the first (and expected only) match so we don't\n  // accidentally touch any other occurrence in the document.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nawait replaceCellValue(\"0.257\", \"0.26\");\nawait replaceCellValue(\"0.513\", \"0.519\");\nawait replaceCellValue(\"0.77\", \"0.779\");\nawait replaceCellValue(\"0.279\", \"0.261\");\n", "ps1": "# Update the p-values in the Mantel correlogram table (x axis label /\n# distance-class rows). Four cells in the \"p\" column change:\n#   0.257 -> 0.26\n#   0.513 -> 0.519\n#   0.77  -> 0.779\n#   0.279 -> 0.261\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText($oldText, $newText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $oldText\n    $range.Find.MatchWholeWord = $true\n    $range.Find.MatchCase = $true\n    $range.Find.Replacement.Text = $newText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $range.Find.Execute(\n        $oldText, $true, $true, $false, $false, $false, $true, 1, $false,\n        $newText, 2\n    )\n}\n\nReplace-ExactText \"0.257\" \"0.26\"\nReplace-ExactText \"0.513\" \"0.519\"\nReplace-ExactText \"0.77\" \"0.779\"\nReplace-ExactText \"0.279\" \"0.261\"\n"}
